# Prozessumbau: Manuelle Zuweisung bei keiner freien Position
#
# This script reproduces the authoring changes made to ManualTests.xlsx:
#  - Three new TODO entries are appended on the "TODO" sheet.
#  - The "Tests" sheet selection moves further down and is no longer the
#    active tab.
#  - The "TODO" sheet becomes the active tab, with an updated selection.

$wb = $excel.ActiveWorkbook

$tests = $wb.Worksheets.Item("Tests")
$todo  = $wb.Worksheets.Item("TODO")

# --- TODO sheet: add the three new rows --------------------------------
# Order of assignment matters because it controls the order new entries
# are appended to the shared string table (A3 -> B1 -> A4).
$todo.Range("A3").Value = "Beim Vorschlagen von Positionen auf das Mindestalter eingehen"
$todo.Range("B1").Value = "OK"
$todo.Range("A4").Value = "Für verfügbare Positionen einen View bauen, der auch die zugehörige Domain enthält/anzeigt"

# --- Tests sheet: update its view / selection (no longer active tab) ---
$tests.Activate() | Out-Null
$tests.Range("A18").Select() | Out-Null

# --- TODO sheet: make it the active tab with new selection --------------
$todo.Activate() | Out-Null
$todo.Range("A5").Select() | Out-Null
